$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
# Column G currently holds "fantasy points" (with the bold header style).
# We need:
#   G1 -> "height"
#   H1 -> "weight"   (new column, needs the same header style as G1)
#   I1 -> "fantasy points" (moved from G1, needs the same header style as G1)

$g1 = $ws.Range("G1")

# Grab the header's formatting before we touch anything, by copying G1's
# format onto the new H1/I1 cells.
$g1.Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value2 = "fantasy points"
$ws.Range("H1").Value2 = "weight"
$ws.Range("G1").Value2 = "height"

# --- Data rows -----------------------------------------------------------
# Column G currently holds each player's "fantasy points" value for that
# game. That value needs to move to column I. Columns G/H then get the
# player's (constant, per-player) height/weight.

$height = 6.416666666666667
$weight = 236

for ($r = 2; $r -le 17; $r++) {
    $fp = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 9).Value2 = $fp
    $ws.Cells.Item($r, 7).Value2 = $height
    $ws.Cells.Item($r, 8).Value2 = $weight
}
